$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column U (21st column)
$ws.Columns("U:U").Insert()

# Set the new header
$ws.Range("U1").Value = "MgCa Coretop modelled temperature"

# Update row 2 values
$ws.Range("R2").Value = 27.73
$ws.Range("S2").Value = -3.698869323730506
$ws.Range("T2").Value = -0.9651359903971048
$ws.Range("U2").Value = 25.1871
$ws.Range("V2").Value = -1.15376667
$ws.Range("W2").Value = 1.579966669999997

